# Remove the two trailing empty paragraphs that were left right after the
# closing "{% endif %}" line / before the section break at the end of the
# table template part. This is part of the "merges tables" speedup work -
# the stray blank paragraphs are no longer needed.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)

if ($last.Range.Text.Trim() -eq "") {
    $last.Range.Delete()
}

$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)

if ($last.Range.Text.Trim() -eq "") {
    $last.Range.Delete()
}
